# Calibration_targets.xlsx update
# - Add CIN3, CIN1, and HPV type-distribution calibration targets (Guan 2012)
#   in rows 177-182 (previously blank placeholder rows).
# - Update the sheet view (scroll position / active selection) to reflect the
#   newly-populated area of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Criteria labels (column A) -- entered grouped by metric (HPV, then CIN1,
# then CIN3; 9v before non-9v within each) so the shared-string table is
# built up in the same order as the source edit.
# ---------------------------------------------------------------------------
$ws.Range("A177").Value = "Proportion HPV attributable to 9v HPV"
$ws.Range("A179").Value = "Proportion CIN1 attributable to 9v HPV"
$ws.Range("A181").Value = "Proportion CIN3 attributable to 9v HPV"
$ws.Range("A178").Value = "Proportion HPV attributable to non-9v HPV"
$ws.Range("A180").Value = "Proportion CIN1 attributable to non-9v HPV"
$ws.Range("A182").Value = "Proportion CIN3 attributable to non-9v HPV"

# Comments (column M) -- entered in row order.
$ws.Range("M177").Value = "Number of women from African studies who were HPV-positive and had normal cytology was 2221."
$ws.Range("M179").Value = "Number of women from African studies who were HPV-positive and had LSIL was 299."
$ws.Range("M181").Value = "Number of women from African studies who were HPV-positive and had HSIL was 185."

# ---------------------------------------------------------------------------
# Source / Group / Year -- same source (Guan 2012 IJC.) and group (all women)
# as the neighbouring rows 175-176, same publication year (2011).
# ---------------------------------------------------------------------------
foreach ($r in 177..182) {
    $ws.Range("B$r").Value = "Guan (2012) IJC."
    $ws.Range("C$r").Value = "all women"
    $ws.Range("D$r").Value = 2011
}

# ---------------------------------------------------------------------------
# Rate / Mean / Variance values.
# ---------------------------------------------------------------------------
$ws.Range("G177").Value = 0.582
$ws.Range("G178").Value = 0.418
$ws.Range("G179").Value = 0.671
$ws.Range("G180").Value = 0.329
$ws.Range("G181").Value = 0.829
$ws.Range("G182").Value = 0.171

$ws.Range("H177").Formula = "=G177"
$ws.Range("H178").Formula = "=G178"
$ws.Range("H179").Formula = "=G179"
$ws.Range("H180").Formula = "=G180"
$ws.Range("H181").Formula = "=G181"
$ws.Range("H182").Formula = "=G182"

$ws.Range("I177").Formula = "=(H177*(1-H177))/2221"
$ws.Range("I178").Formula = "=(H178*(1-H178))/2221"
$ws.Range("I179").Formula = "=(H179*(1-H179))/299"
$ws.Range("I180").Formula = "=(H180*(1-H180))/299"
$ws.Range("I181").Formula = "=(H181*(1-H181))/185"
$ws.Range("I182").Formula = "=(H182*(1-H182))/185"

# Wrap the "Comments" column cells (M176:M182), matching the formatting of
# the rest of the Comments column.
$ws.Range("M176:M182").WrapText = $true

# Row heights: 177 & 179 wrap onto two lines (2 x 14.5 = 29pt); 181 & 182 are
# explicitly sized.
$ws.Rows.Item(177).RowHeight = 29
$ws.Rows.Item(179).RowHeight = 29
$ws.Rows.Item(181).RowHeight = 29.5
$ws.Rows.Item(182).RowHeight = 15

# ---------------------------------------------------------------------------
# Sheet view: scrolled down to the newly-added rows, selection moved to G189
# ---------------------------------------------------------------------------
$ws.Range("G189").Select() | Out-Null

$wb.Save()
